# Add a "Save" column (H) to the s_vals sheet, matching the style of the
# existing header row (column G = "sum").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1) onto the new
# header cell (H1) so it reuses the same cell style (bold, bordered,
# centered) instead of Excel minting a brand-new style entry.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data values for the Save column.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
